$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 土地 (land) - used as a source of already-correct header/common-field cells
$ws3 = $wb.Worksheets.Item(3)   # 汽車 (car) - the sheet being fixed

# -----------------------------------------------------------------
# Row 1 (header row) on the 汽車 sheet currently (incorrectly) holds
# copied data values instead of field-name labels. Fix it up and
# extend it with the common trailing columns (H1:N1) that already
# exist, correctly, on the 土地 sheet.
# -----------------------------------------------------------------

# B1:G1 -> name, capacity(new), owner, register_date, register_reason, acquire_value
$ws1.Range("B1").Copy($ws3.Range("B1"))
$ws3.Range("C1").Value = "capacity"
$ws1.Range("E1").Copy($ws3.Range("D1"))
$ws1.Range("F1").Copy($ws3.Range("E1"))
$ws1.Range("G1").Copy($ws3.Range("F1"))
$ws1.Range("H1").Copy($ws3.Range("G1"))

# H1:N1 -> property_category, category, date, legislator_name, legislator_id, source_file, index
$ws1.Range("I1").Copy($ws3.Range("H1"))
$ws1.Range("J1").Copy($ws3.Range("I1"))
$ws1.Range("K1").Copy($ws3.Range("J1"))
$ws1.Range("L1").Copy($ws3.Range("K1"))
$ws1.Range("M1").Copy($ws3.Range("L1"))
$ws1.Range("N1").Copy($ws3.Range("M1"))
$ws1.Range("O1").Copy($ws3.Range("N1"))

# -----------------------------------------------------------------
# Row 2 (data row): keep name/owner/register_date/register_reason/
# acquire_value as they are; turn the "2987" engine displacement
# text into a real number under the new capacity column, and append
# the common trailing fields (H2:N2).
# -----------------------------------------------------------------

$ws3.Range("C2").Value = 2987

$ws1.Range("I2").Copy($ws3.Range("H2"))
$ws1.Range("J2").Copy($ws3.Range("I2"))
$ws1.Range("K2").Copy($ws3.Range("J2"))
$ws1.Range("L2").Copy($ws3.Range("K2"))
$ws3.Range("L2").Value = 629
$ws1.Range("N2").Copy($ws3.Range("M2"))
$ws3.Range("N2").Value = 27
